# Fruta / hortaliza, semanal
#
# Weekly refresh: a new price observation (2021-11-30) is inserted at row 11,
# pushing the existing historical rows (11-119) down by one row (to 12-120).
# The brand-new row 11 keeps the same categorical values (market, region,
# product, quality, volume, unit, origin, kg/unit) that the old row 11 had,
# but carries fresh date / price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 11:119 down to 12:120 (whole-row insert, like Excel's
# "Insert Sheet Rows" command) so every column -- including the ones that
# don't change week to week -- moves together.
$ws.Rows(11).Insert()

# Populate the newly-inserted row 11 with this week's record.
$ws.Range("A11").Value = 11
$ws.Range("B11").Value = 'Vega Monumental Concepción'
$ws.Range("C11").Value = 'Bíobío'
$ws.Range("D11").Value = 44530
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 'Fruta'
$ws.Range("G11").Value = 100108
$ws.Range("H11").Value = 'Tropicales y subtropicales'
$ws.Range("I11").Value = 100108005
$ws.Range("J11").Value = 'Piña'
$ws.Range("K11").Value = 'Caramelo'
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 200
$ws.Range("N11").Value = 18500
$ws.Range("O11").Value = 19000
$ws.Range("P11").Value = 18750
$ws.Range("Q11").Value = '$/caja 14 unidades'
$ws.Range("R11").Value = 'Ecuador'
$ws.Range("S11").Value = 1339
$ws.Range("T11").Value = 14
